# Auto-generated Excel COM-interop script to apply scheduled-runner price updates
# across the Zodiark_Profits workbook's 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value2 = 2117.516
$ws.Range("I15").Value2 = 2117.516
$ws.Range("K15").Value2 = 6352.548000000001
$ws.Range("M15").Value2 = -6183.548000000001
# Row 43
$ws.Range("H43").Value2 = 1691
$ws.Range("J43").Value2 = 2624.5
$ws.Range("L43").Value2 = 2624.5
$ws.Range("N43").Value2 = -2762.5
# Row 49
$ws.Range("H49").Value2 = 19
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 19
$ws.Range("K49").Value2 = 0
$ws.Range("L49").Value2 = 57
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value2 = -329
# Row 52
$ws.Range("H52").Value2 = 6097.75
$ws.Range("I52").Value2 = 5997.5
$ws.Range("J52").Value2 = 6198
$ws.Range("K52").Value2 = 17992.5
$ws.Range("L52").Value2 = 18594
$ws.Range("M52").Value2 = -17832.5
$ws.Range("N52").Value2 = -18914
# Row 59
$ws.Range("H59").Value2 = 387
$ws.Range("I59").Value2 = 530.5
$ws.Range("J59").Value2 = 100
$ws.Range("K59").Value2 = 1591.5
$ws.Range("L59").Value2 = 300
$ws.Range("M59").Value2 = -1034.5
$ws.Range("N59").Value2 = -1414
# Row 63
$ws.Range("H63").Value2 = 100000
$ws.Range("J63").Value2 = 100000
$ws.Range("L63").Value2 = 100000
$ws.Range("N63").Value2 = -101248
# Row 66
$ws.Range("H66").Value2 = 100000
$ws.Range("J66").Value2 = 100000
$ws.Range("L66").Value2 = 300000
$ws.Range("N66").Value2 = -306240
# Row 98
$ws.Range("H98").Value2 = 2005.375
$ws.Range("I98").Value2 = 2027.3478
$ws.Range("K98").Value2 = 2027.3478
$ws.Range("M98").Value2 = -529.3478
# Row 116
$ws.Range("H116").Value2 = 3410.7585
$ws.Range("I116").Value2 = 3409
$ws.Range("J116").Value2 = 3426
$ws.Range("K116").Value2 = 3409
$ws.Range("L116").Value2 = 3426
$ws.Range("M116").Value2 = 33
$ws.Range("N116").Value2 = -10310
# Row 122
$ws.Range("H122").Value2 = 2005.375
$ws.Range("I122").Value2 = 2027.3478
$ws.Range("K122").Value2 = 6082.0434
$ws.Range("M122").Value2 = -3632.0434
# Row 132
$ws.Range("H132").Value2 = 19609068
$ws.Range("I132").Value2 = 19609068
$ws.Range("K132").Value2 = 58827204
$ws.Range("M132").Value2 = -58824674

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 4030.3171
$ws.Range("I32").Value2 = 3187.0303
$ws.Range("K32").Value2 = 3187.0303
$ws.Range("M32").Value2 = -2900.0303
# Row 45
$ws.Range("H45").Value2 = 1206.8096
$ws.Range("J45").Value2 = 1102.3334
$ws.Range("L45").Value2 = 1102.3334
$ws.Range("N45").Value2 = -1856.3334
# Row 74
$ws.Range("H74").Value2 = 2781.9678
$ws.Range("I74").Value2 = 2491.875
$ws.Range("K74").Value2 = 2491.875
$ws.Range("M74").Value2 = -1617.875
# Row 77
$ws.Range("H77").Value2 = 2781.9678
$ws.Range("I77").Value2 = 2491.875
$ws.Range("K77").Value2 = 12459.375
$ws.Range("M77").Value2 = -8091.375
# Row 97
$ws.Range("H97").Value2 = 582.1667
$ws.Range("I97").Value2 = 582.1667
$ws.Range("K97").Value2 = 582.1667
$ws.Range("M97").Value2 = -86.16669999999999
# Row 102
$ws.Range("H102").Value2 = 27807624
$ws.Range("I102").Value2 = 33335688
$ws.Range("K102").Value2 = 33335688
$ws.Range("M102").Value2 = -33334066

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value2 = 71431360
$ws.Range("I20").Value2 = 111113550
$ws.Range("K20").Value2 = 111113550
$ws.Range("M20").Value2 = -111113303

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value2 = 862.13043
$ws.Range("J22").Value2 = 1636.1111
$ws.Range("L22").Value2 = 1636.1111
$ws.Range("N22").Value2 = -2336.1111
# Row 31
$ws.Range("H31").Value2 = 2267.0605
$ws.Range("I31").Value2 = 2171.087
$ws.Range("J31").Value2 = 2487.8
$ws.Range("K31").Value2 = 2171.087
$ws.Range("L31").Value2 = 2487.8
$ws.Range("M31").Value2 = -1876.087
$ws.Range("N31").Value2 = -3077.8
# Row 34
$ws.Range("H34").Value2 = 2267.0605
$ws.Range("I34").Value2 = 2171.087
$ws.Range("J34").Value2 = 2487.8
$ws.Range("K34").Value2 = 2171.087
$ws.Range("L34").Value2 = 2487.8
$ws.Range("M34").Value2 = -1969.087
$ws.Range("N34").Value2 = -2891.8
# Row 39
$ws.Range("H39").Value2 = 14500
$ws.Range("I39").Value2 = 14500
$ws.Range("K39").Value2 = 14500
$ws.Range("M39").Value2 = -14109
# Row 49
$ws.Range("H49").Value2 = 14500
$ws.Range("I49").Value2 = 14500
$ws.Range("K49").Value2 = 14500
$ws.Range("M49").Value2 = -14318
# Row 69
$ws.Range("H69").Value2 = 37920.375
$ws.Range("I69").Value2 = 49872.8
$ws.Range("K69").Value2 = 49872.8
$ws.Range("M69").Value2 = -49123.8
# Row 72
$ws.Range("H72").Value2 = 37920.375
$ws.Range("I72").Value2 = 49872.8
$ws.Range("K72").Value2 = 149618.4
$ws.Range("M72").Value2 = -145874.4
# Row 99
$ws.Range("H99").Value2 = 4161.4
$ws.Range("I99").Value2 = 2269.3333
$ws.Range("K99").Value2 = 2269.3333
$ws.Range("M99").Value2 = -771.3332999999998
# Row 122
$ws.Range("H122").Value2 = 3023.25
$ws.Range("I122").Value2 = 1681.7273
$ws.Range("K122").Value2 = 5045.1819
$ws.Range("M122").Value2 = -2595.1819
# Row 126
$ws.Range("H126").Value2 = 4161.4
$ws.Range("I126").Value2 = 2269.3333
$ws.Range("K126").Value2 = 6807.999899999999
$ws.Range("M126").Value2 = -4337.999899999999

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Range("H94").Value2 = 0
$ws.Range("J94").Value2 = 0
$ws.Range("L94").Value2 = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value2 = 40076
$ws.Range("I70").Value2 = 111186.6
$ws.Range("J70").Value2 = 12725.77
$ws.Range("K70").Value2 = 111186.6
$ws.Range("L70").Value2 = 12725.77
$ws.Range("M70").Value2 = -110916.6
$ws.Range("N70").Value2 = -13265.77
# Row 73
$ws.Range("H73").Value2 = 40076
$ws.Range("I73").Value2 = 111186.6
$ws.Range("J73").Value2 = 12725.77
$ws.Range("K73").Value2 = 111186.6
$ws.Range("L73").Value2 = 12725.77
$ws.Range("M73").Value2 = -110250.6
$ws.Range("N73").Value2 = -14597.77
# Row 107
$ws.Range("H107").Value2 = 1636.1538
$ws.Range("I107").Value2 = 1390
$ws.Range("K107").Value2 = 1390
$ws.Range("M107").Value2 = 530
# Row 126
$ws.Range("H126").Value2 = 5532.2354
$ws.Range("I126").Value2 = 4388.923
$ws.Range("K126").Value2 = 13166.769
$ws.Range("M126").Value2 = -10696.769

$ws = $wb.Worksheets.Item("LTW")
# Row 76
$ws.Range("H76").Value2 = 29999.5
$ws.Range("J76").Value2 = 29999.5
$ws.Range("L76").Value2 = 29999.5
$ws.Range("N76").Value2 = -30675.5
# Row 79
$ws.Range("H79").Value2 = 29999.5
$ws.Range("J79").Value2 = 29999.5
$ws.Range("L79").Value2 = 29999.5
$ws.Range("N79").Value2 = -32339.5
# Row 101
$ws.Range("H101").Value2 = 91323.14
$ws.Range("J101").Value2 = 91323.14
$ws.Range("L101").Value2 = 91323.14
$ws.Range("N101").Value2 = -97813.14
# Row 104
$ws.Range("H104").Value2 = 16677.428
$ws.Range("J104").Value2 = 16677.428
$ws.Range("L104").Value2 = 16677.428
$ws.Range("N104").Value2 = -23665.428

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value2 = 3095.5
$ws.Range("I81").Value2 = 795
$ws.Range("J81").Value2 = 5396
$ws.Range("K81").Value2 = 1590
$ws.Range("L81").Value2 = 10792
$ws.Range("M81").Value2 = -529
$ws.Range("N81").Value2 = -12914
# Row 84
$ws.Range("H84").Value2 = 3095.5
$ws.Range("I84").Value2 = 795
$ws.Range("J84").Value2 = 5396
$ws.Range("K84").Value2 = 7950
$ws.Range("L84").Value2 = 53960
$ws.Range("M84").Value2 = -2646
$ws.Range("N84").Value2 = -64568
# Row 132
$ws.Range("H132").Value2 = 1590.8462
$ws.Range("I132").Value2 = 1598.5
$ws.Range("K132").Value2 = 4795.5
$ws.Range("M132").Value2 = -2265.5
# Row 136
$ws.Range("H136").Value2 = 3614.1875
$ws.Range("I136").Value2 = 3065.3809
$ws.Range("K136").Value2 = 9196.1427
$ws.Range("M136").Value2 = -6646.1427
